$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update price values as described in the diff
$ws.Range("C2").Value = 9099
$ws.Range("C3").Value = 12299

# Update the active selection to C3 (matches the sheetView selection change)
$ws.Range("C3").Select()
